# bug fixes, added del and rename
# Adds two new keyword rows (DEL / RENAME) to the "Keywords" sheet, mirroring
# the existing "IF" row's layout, and fills in the previously-missing
# Parse/Eval/String/Destroy ("X") markers on the "IF" row itself.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Keywords")

# --- Row 19 (IF): the Parse/Eval/String/Destroy columns (C:F) were blank; mark them "X" ---
$ws.Range("C19").Value = "X"
$ws.Range("D19").Value = "X"
$ws.Range("E19").Value = "X"
$ws.Range("F19").Value = "X"

# --- Row 20 (new): DEL keyword ---
$ws.Range("B20").Value = "DEL"
$ws.Range("C20").Value = "X"
$ws.Range("D20").Value = "X"
$ws.Range("E20").Value = "X"
$ws.Range("F20").Value = "X"
$ws.Range("G20").Value = "BTOKEN_DEL"
$ws.Range("H20").Value = "EXPRINDEX"
$ws.Range("H20:K20").Merge()
$ws.Range("H20:K20").HorizontalAlignment = -4108

# --- Row 21 (new): RENAME keyword ---
$ws.Range("B21").Value = "RENAME"
$ws.Range("C21").Value = "X"
$ws.Range("D21").Value = "X"
$ws.Range("E21").Value = "X"
$ws.Range("F21").Value = "X"
$ws.Range("G21").Value = "BTOKEN_RENAME"

$ws.Range("H21").Value = "EXPRINDEX"
$ws.Range("H21:K21").Merge()
$ws.Range("H21:K21").HorizontalAlignment = -4108

$ws.Range("L21").Value = "EXPRINDEX"
$ws.Range("L21:O21").Merge()
$ws.Range("L21:O21").HorizontalAlignment = -4108

# --- Update the saved selection/active cell to reflect where the edit left off ---
$ws.Activate()
$ws.Range("P21").Select()
